# added report write to clearCartItem()
# Flip the Execution Flag (column C) results for the test-suite rows so the
# newly added report-write call in clearCartItem() gets exercised:
#   Search Page (row4), COD Order (row6), Checking Filters (row7) and
#   Emailing Reports (row8) now run ("YES"); Product Listing (row5) is
#   turned off ("NO").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = "YES"
$ws.Range("C5").Value = "NO"
$ws.Range("C6").Value = "YES"
$ws.Range("C7").Value = "YES"
$ws.Range("C8").Value = "YES"

# Move the selection/cursor as recorded in the saved view state.
$ws.Range("E15").Select()
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 1
